$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2023" data column (M) by mirroring column L's styling ---
# Row 3 (empty header separator cell)
$ws.Range("L3").Copy($ws.Range("M3"))

# Row 4 (year header)
$ws.Range("L4").Copy($ws.Range("M4"))
$ws.Range("M4").Value = 2023

# Row 5 (Total)
$ws.Range("L5").Copy($ws.Range("M5"))
$ws.Range("M5").Value = 311.65582791395695

# Row 6 (including: -- header, empty data cell)
$ws.Range("L6").Copy($ws.Range("M6"))

# Row 7 (solid)
$ws.Range("L7").Copy($ws.Range("M7"))
$ws.Range("M7").Value = 119.55977988994496

# Row 8 (gaseous and liquid)
$ws.Range("L8").Copy($ws.Range("M8"))
$ws.Range("M8").Value = 192.09604802401199

# Row 9 (of them: -- header, empty data cell)
$ws.Range("L9").Copy($ws.Range("M9"))

# Row 10 (sulfur dioxide)
$ws.Range("L10").Copy($ws.Range("M10"))
$ws.Range("M10").Value = 78.539269634817401

# Row 11 (carbon monoxide)
$ws.Range("L11").Copy($ws.Range("M11"))
$ws.Range("M11").Value = 60.030015007503756

# Row 12 (nitrogen oxides)
$ws.Range("L12").Copy($ws.Range("M12"))
$ws.Range("M12").Value = 26.013006503251628

# --- Row heights: rows 2, 5-11 pick up an explicit 15pt height; row 12's
#     explicit height grows from 13.5 to 15 ---
$ws.Range("A2").RowHeight = 15
$ws.Range("A5").RowHeight = 15
$ws.Range("A6").RowHeight = 15
$ws.Range("A7").RowHeight = 15
$ws.Range("A8").RowHeight = 15
$ws.Range("A9").RowHeight = 15
$ws.Range("A10").RowHeight = 15
$ws.Range("A11").RowHeight = 15
$ws.Range("A12").RowHeight = 15

# --- Column A:C width shrinks slightly ---
$ws.Range("A1:C1").ColumnWidth = 39.6

# --- Clear the lingering N5 selection left over from the previous save ---
$ws.Range("A1").Select()
